$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new D value (or $null to leave D unchanged), new E value (without padding spaces)
$updates = @(
    @{row=2; d="61.027.17"; e="+0.45%"},
    @{row=3; d="3.387.63"; e="+0.17%"},
    @{row=4; d=$null; e="-0.05%"},
    @{row=5; d="572.03"; e="+0.11%"},
    @{row=6; d="141.87"; e="+0.37%"},
    @{row=7; d=$null; e="-0.01%"},
    @{row=8; d="0.475"; e="+0.18%"},
    @{row=9; d=$null; e="+2.07%"},
    @{row=10; d=$null; e="-0.83%"},
    @{row=11; d=$null; e="-1.14%"},
    @{row=12; d="3.965.96"; e="+0.13%"},
    @{row=13; d=$null; e="+1.91%"},
    @{row=14; d=$null; e="-0.81%"},
    @{row=15; d=$null; e="+0.29%"},
    @{row=16; d="3.379.23"; e="-0.26%"},
    @{row=17; d="61.104.14"; e="+0.34%"},
    @{row=18; d="6.11"; e="-2.59%"},
    @{row=19; d="13.67"; e="-3.11%"},
    @{row=20; d="8.91"; e="-1.30%"},
    @{row=21; d="383.05"; e="-1.40%"},
    @{row=22; d="75.29"; e="+2.60%"},
    @{row=23; d=$null; e="-1.31%"},
    @{row=24; d=$null; e="-0.05%"},
    @{row=25; d=$null; e="-1.20%"},
    @{row=26; d="3.522.90"; e="+0.01%"},
    @{row=27; d=$null; e="+2.50%"},
    @{row=28; d=$null; e="+0.04%"},
    @{row=29; d="7.24"; e="-2.28%"},
    @{row=30; d="7.97"; e="-1.25%"},
    @{row=31; d="2.15"; e="-1.14%"},
    @{row=32; d=$null; e="-0.03%"},
    @{row=33; d="1.38"; e="-4.27%"},
    @{row=34; d="23.21"; e="-2.09%"},
    @{row=35; d=$null; e="+0.38%"},
    @{row=36; d="166.62"; e="-0.31%"},
    @{row=37; d="3.419.74"; e="+0.21%"},
    @{row=38; d="4.98"; e="-0.78%"},
    @{row=39; d=$null; e="-2.54%"},
    @{row=40; d="0.0767"; e="-1.18%"},
    @{row=41; d="26.77"; e="-0.66%"},
    @{row=42; d=$null; e="-0.06%"},
    @{row=43; d="0.779"; e="-0.40%"},
    @{row=44; d=$null; e="-1.84%"},
    @{row=45; d=$null; e="-1.55%"},
    @{row=46; d="1.13"; e="+0.13%"},
    @{row=47; d="2.449.47"; e="-3.36%"},
    @{row=48; d="22.97"; e="+0.30%"},
    @{row=49; d="6.69"; e="-1.91%"},
    @{row=50; d="2.16"; e="+10.96%"},
    @{row=51; d=$null; e="+0.14%"}
)

foreach ($item in $updates) {
    if ($item.d -ne $null) {
        $dc = $ws.Range("D$($item.row)")
        $dc.Value = "'" + $item.d
        $dc.Style = "Normal"
    }
    $ec = $ws.Range("E$($item.row)")
    $ec.Value = "'" + "  " + $item.e + "  "
    $ec.Style = "Normal"
}